$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.49"
$ws.Range("E2").Value = "'-4.27%"
$ws.Range("D3").Value = "'30.40"
$ws.Range("E3").Value = "'-6.42%"
$ws.Range("D4").Value = "'4.941"
$ws.Range("E4").Value = "'-2.48%"
$ws.Range("D5").Value = "'0.07180"
$ws.Range("E5").Value = "'-6.54%"
$ws.Range("D6").Value = "'1.794"
$ws.Range("E6").Value = "'-14.51%"
$ws.Range("D7").Value = "'7.598"
$ws.Range("E7").Value = "'-3.12%"
$ws.Range("D8").Value = "'3.716"
$ws.Range("E8").Value = "'-1.65%"
$ws.Range("D9").Value = "'0.8992"
$ws.Range("E9").Value = "'-2.29%"
$ws.Range("D10").Value = "'0.1659"
$ws.Range("E10").Value = "'-6.16%"
$ws.Range("D11").Value = "'0.07767"
$ws.Range("E11").Value = "'-1.63%"
$ws.Range("D12").Value = "'0.08005"
$ws.Range("E12").Value = "'-5.15%"
$ws.Range("D13").Value = "'0.03035"
$ws.Range("E13").Value = "'-0.78%"
$ws.Range("D14").Value = "'0.09993"
$ws.Range("E14").Value = "'0.07%"
$ws.Range("D15").Value = "'0.001507"
$ws.Range("E15").Value = "'-0.58%"
$ws.Range("B16").Value = "'CoinExToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04502"
$ws.Range("E16").Value = "'-0.71%"
$ws.Range("B17").Value = "'TigerCash"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005838"
$ws.Range("E17").Value = "'-0.49%"
$ws.Range("B18").Value = "'UpBots"
$ws.Range("C18").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D18").Value = "'0.007492"
$ws.Range("E18").Value = "'-0.08%"
$ws.Range("B19").Value = "'LEO"
$ws.Range("C19").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.479"
$ws.Range("E19").Value = "'0.48%"
$ws.Range("B20").Value = "'BTSEToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.067"
$ws.Range("E20").Value = "'-3.90%"
$ws.Range("B21").Value = "'BitpandaEcosystemToken"
$ws.Range("C21").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3302"
$ws.Range("E21").Value = "'-1.16%"
$ws.Range("B22").Value = "'ProBitToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1295"
$ws.Range("E22").Value = "'-2.53%"
$ws.Range("B23").Value = "'MCDex"
$ws.Range("C23").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").Value = "'3.962"
$ws.Range("E23").Value = "'-7.25%"
$ws.Range("B24").Value = "'ZBToken"
$ws.Range("C24").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.2103"
$ws.Range("E24").Value = "'6.39%"
$ws.Range("D25").Value = "'0.001214"
$ws.Range("E25").Value = "'-1.89%"
$ws.Range("D26").Value = "'0.004620"
$ws.Range("E26").Value = "'11.99%"
$ws.Range("D27").Value = "'0.0001299"
$ws.Range("E27").Value = "'3.90%"
$ws.Range("D39").Value = "'0.01551"
$ws.Range("E39").Value = "'-9.31%"
$ws.Range("D40").Value = "'0.04310"
$ws.Range("E40").Value = "'-7.77%"
$ws.Range("D41").Value = "'0.007320"
$ws.Range("E41").Value = "'-2.04%"
$ws.Range("D43").Value = "'0.1298"
$ws.Range("E43").Value = "'-4.13%"
$ws.Range("D44").Value = "'0.002012"
$ws.Range("E44").Value = "'-13.70%"
$ws.Range("D45").Value = "'0.009104"
$ws.Range("E45").Value = "'-13.52%"
$ws.Range("D46").Value = "'0.00005895"
$ws.Range("E46").Value = "'-5.04%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'2.255"
$ws.Range("E48").Value = "'308.14%"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.06%"
